$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.596.33'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '2.011.59'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.32'
$ws.Range("E5").Value = '  -1.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.633'
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '62.69'
$ws.Range("E7").Value = '  +2.16%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +4.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.15'
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0787'
$ws.Range("E11").Value = '  +6.33%  '
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.882'
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.96'
$ws.Range("E14").Value = '  +13.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.16'
$ws.Range("E15").Value = '  -4.36%  '
$ws.Range("D16").Value = '2.304.24'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.53'
$ws.Range("E17").Value = '  +1.99%  '
$ws.Range("D18").Value = '2.011.36'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").Value = '36.529.65'
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.92'
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").Value = '0.0₃0875'
$ws.Range("E21").Value = '  +1.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.34'
$ws.Range("E22").Value = '  +2.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.97'
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.53'
$ws.Range("E25").Value = '  -8.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.99'
$ws.Range("E27").Value = '  +4.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.36'
$ws.Range("E28").Value = '  -2.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.16'
$ws.Range("E29").Value = '  +2.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.132'
$ws.Range("E30").Value = '  +21.18%  '
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("E32").Value = '  -1.51%  '
$ws.Range("E33").Value = '  -0.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0633'
$ws.Range("E34").Value = '  +4.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.51'
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.47'
$ws.Range("E36").Value = '  +9.77%  '
$ws.Range("E37").Value = '  -3.78%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.18'
$ws.Range("E40").Value = '  +13.68%  '
$ws.Range("E41").Value = '  -2.06%  '
$ws.Range("E42").Value = '  +3.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.90'
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0217'
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '95.92'
$ws.Range("E46").Value = '  +1.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.72'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.69'
$ws.Range("E48").Value = '  -4.56%  '
$ws.Range("D49").Value = '1.361.93'
$ws.Range("E49").Value = '  -5.22%  '
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("D51").Value = '2.196.71'
$ws.Range("E51").Value = '  +0.25%  '
